$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Semestre ideal" value EQD-8,EQN-11 -> EQD-8,EQN-10
$ws.Range("B9").Value = "EQD-8,EQN-10"
$ws.Range("C9").Value = "EQD-8,EQN-10"

# Update the last requisito row (row 25) text
$ws.Range("B25").Value = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n"

# Add new row 26 with new requisito text, matching row 25's formatting/style
$ws.Range("B26").Value = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"
$ws.Range("C26").Value = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"

$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("C25").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Rows.Item(26).RowHeight = $ws.Rows.Item(25).RowHeight
